{"js": "// The \"Normal Text\" example paragraph was split across four runs\n// (an artifact of earlier edits). Word renders it seamlessly as one\n// sentence, but the underlying OOXML carries unnecessary run splits.\n// Collapse it back down to a single run without changing the visible\n// text, so the template's styles render/apply cleanly.\nconst targetText =\n  \"This is normal, editable text that reviewers should feel comfortable editing.\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text === targetText);\n\nif (target) {\n  // Replacing a paragraph's range with its own text forces Word to\n  // rewrite it as a single run, merging the fragmented runs.\n  target.insertText(targetText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The \"Normal Text\" example paragraph was split across four runs\n# (an artifact of earlier edits). Word renders it seamlessly as one\n# sentence, but the underlying OOXML carries unnecessary run splits.\n# Collapse it back down to a single run without changing the visible\n# text, so the template's styles render/apply cleanly.\n\n$targetText = \"This is normal, editable text that reviewers should feel comfortable editing.\"\n\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $paraText = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($paraText -eq $targetText) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # Running Find/Replace over the paragraph's own text forces Word to\n    # rewrite it as a single run, merging the fragmented runs even\n    # though the visible text is unchanged.\n    $find = $target.Range.Find\n    $find.Execute($targetText, $false, $false, $false, $false, $false, $true, 1, $false, $targetText, 2)\n}\n"}
